# Add "Sentiment Analysis" as a new column to the Table2 table (expands
# A1:Q16 -> A1:R16) and populate it with the sentiment-analysis values,
# matching the "added sentiment analysis to report" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one column; Excel/the table engine takes care of
# updating the table ref, autoFilter ref, dimension, etc.
$newCol = $lo.ListColumns.Add()
$newCol.Range.Item(1).Value = "Sentiment Analysis"

# Sentiment Analysis values, row by row (rows 2-16, column R / 18)
$values = @(
  0.0918036076472,
  0.00380666456608,
  0.112327180368,
  0.0961059135493,
  0.0992528586376,
  0.0924029421256,
  0.148926322657,
  0.118894264242,
  0.0810604753697,
  0.148720333191,
  0.169905447166,
  0.0940876494875,
  0.112327180368,
  0.0747258903217,
  0.138129282178
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $c = $ws.Cells.Item($row, 18)
    $c.Value = $values[$i]
    $c.HorizontalAlignment = -4108   # xlCenter, matches the other data columns
}

# Match the column width Excel picked for the new "Sentiment Analysis" column
$ws.Columns.Item(18).ColumnWidth = 30.33

# Keep the new cell selected, like the saved file shows
$ws.Range("R16").Select()
